$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($ws, $addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-CellText $ws "D2" "68.968.86"
Set-CellText $ws "E2" "  -5.44%  "

Set-CellText $ws "D3" "3.744.03"
Set-CellText $ws "E3" "  -6.03%  "

Set-CellText $ws "E4" "  +0.53%  "

Set-CellText $ws "D5" "580.88"
Set-CellText $ws "E5" "  -6.03%  "

Set-CellText $ws "D6" "177.66"
Set-CellText $ws "E6" "  +7.57%  "

Set-CellText $ws "D7" "0.645"
Set-CellText $ws "E7" "  -5.91%  "

Set-CellText $ws "D8" "1.01"
Set-CellText $ws "E8" "  +0.68%  "

Set-CellText $ws "D9" "0.728"
Set-CellText $ws "E9" "  -4.18%  "

Set-CellText $ws "D10" "0.170"
Set-CellText $ws "E10" "  +1.07%  "

Set-CellText $ws "D11" "53.32"
Set-CellText $ws "E11" "  -9.12%  "

Set-CellText $ws "D12" "0.0000307"
Set-CellText $ws "E12" "  -3.20%  "

Set-CellText $ws "D13" "10.90"
Set-CellText $ws "E13" "  -2.73%  "

Set-CellText $ws "D14" "4.359.62"
Set-CellText $ws "E14" "  -5.77%  "

Set-CellText $ws "D15" "3.806.81"
Set-CellText $ws "E15" "  -4.84%  "

Set-CellText $ws "D16" "19.71"
Set-CellText $ws "E16" "  -4.53%  "

Set-CellText $ws "D17" "13.24"
Set-CellText $ws "E17" "  -7.59%  "

Set-CellText $ws "D18" "1.16"
Set-CellText $ws "E18" "  -8.77%  "

Set-CellText $ws "E19" "  -2.72%  "

Set-CellText $ws "D20" "69.176.72"
Set-CellText $ws "E20" "  -5.08%  "

Set-CellText $ws "D21" "417.06"
Set-CellText $ws "E21" "  -5.33%  "

Set-CellText $ws "E22" "  -6.70%  "

Set-CellText $ws "D23" "90.40"
Set-CellText $ws "E23" "  -6.17%  "

Set-CellText $ws "D24" "3.12"
Set-CellText $ws "E24" "  -8.22%  "

Set-CellText $ws "D25" "13.15"
Set-CellText $ws "E25" "  -8.66%  "

Set-CellText $ws "D26" "11.04"
Set-CellText $ws "E26" "  -2.34%  "

Set-CellText $ws "D27" "3.82"
Set-CellText $ws "E27" "  -7.74%  "

Set-CellText $ws "D28" "5.94"
Set-CellText $ws "E28" "  -0.41%  "

Set-CellText $ws "D29" "9.90"
Set-CellText $ws "E29" "  -5.99%  "

Set-CellText $ws "D30" "8.41"
Set-CellText $ws "E30" "  +7.09%  "

Set-CellText $ws "D31" "33.57"
Set-CellText $ws "E31" "  -7.23%  "

Set-CellText $ws "D32" "13.04"
Set-CellText $ws "E32" "  -5.89%  "

Set-CellText $ws "D33" "45.63"
Set-CellText $ws "E33" "  -7.25%  "

Set-CellText $ws "D34" "0.120"
Set-CellText $ws "E34" "  -8.12%  "

Set-CellText $ws "D35" "66.66"
Set-CellText $ws "E35" "  -6.86%  "

Set-CellText $ws "D36" "0.0₃0946"
Set-CellText $ws "E36" "  +1.90%  "

Set-CellText $ws "D37" "616.51"
Set-CellText $ws "E37" "  -3.95%  "

Set-CellText $ws "D38" "0.409"
Set-CellText $ws "E38" "  -6.05%  "

Set-CellText $ws "D39" "0.998"
Set-CellText $ws "E39" "  -0.13%  "

Set-CellText $ws "D40" "1.01"
Set-CellText $ws "E40" "  +0.43%  "

Set-CellText $ws "B41" "Kaspa"
Set-CellText $ws "C41" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-CellText $ws "D41" "0.140"
Set-CellText $ws "E41" "  -5.10%  "

Set-CellText $ws "B42" "dogwifhat"
Set-CellText $ws "C42" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-CellText $ws "D42" "3.21"
Set-CellText $ws "E42" "  +9.20%  "

Set-CellText $ws "D43" "3.14"
Set-CellText $ws "E43" "  -10.88%  "

Set-CellText $ws "B44" "EnergySwap"
Set-CellText $ws "C44" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-CellText $ws "D44" "30.13"
Set-CellText $ws "E44" "  +18.04%  "

Set-CellText $ws "B45" "VeChain"
Set-CellText $ws "C45" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-CellText $ws "D45" "0.0449"
Set-CellText $ws "E45" "  -7.95%  "

Set-CellText $ws "B46" "THORChain"
Set-CellText $ws "C46" "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-CellText $ws "D46" "9.71"
Set-CellText $ws "E46" "  -12.04%  "

Set-CellText $ws "B47" "Fetch.AI"
Set-CellText $ws "C47" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-CellText $ws "D47" "2.62"
Set-CellText $ws "E47" "  -0.48%  "

Set-CellText $ws "B48" "Stellar"
Set-CellText $ws "C48" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-CellText $ws "D48" "0.138"
Set-CellText $ws "E48" "  -7.53%  "

Set-CellText $ws "D49" "3.24"
Set-CellText $ws "E49" "  -5.16%  "

Set-CellText $ws "B50" "WEMIXToken"
Set-CellText $ws "C50" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-CellText $ws "D50" "2.73"
Set-CellText $ws "E50" "  -16.63%  "

Set-CellText $ws "B51" "Maker"
Set-CellText $ws "C51" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-CellText $ws "D51" "2.778.29"
Set-CellText $ws "E51" "  -4.54%  "
